$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = [double]"8.635273788982842e-21"
$ws.Range("C2").Value = [double]"1.609528544852453e-10"
$ws.Range("D2").Value = 1

$ws.Range("A3").Value = -1
$ws.Range("B3").Value = [double]"-2.414292817278677e-10"
$ws.Range("C3").Value = [double]"5.233679649020949e-17"
$ws.Range("D3").Value = [double]"-7.428594357819736e-11"

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = [double]"3.219057089704903e-10"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = [double]"-1.609528544852453e-10"

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = [double]"-3.219057089704904e-10"
$ws.Range("D5").Value = [double]"2.524156694435312e-20"
